$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.701.01'
$ws.Range('E2').Value = '  -2.60%  '
$ws.Range('D3').Value = '3.555.29'
$ws.Range('E3').Value = '  -3.29%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '616.86'
$ws.Range('E5').Value = '  -5.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.26'
$ws.Range('E6').Value = '  -3.32%  '
$ws.Range('D7').Value = '3.552.54'
$ws.Range('E7').Value = '  -3.17%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.485'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.89'
$ws.Range('E11').Value = '  -2.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = '4.152.94'
$ws.Range('E14').Value = '  -3.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.18'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '3.548.32'
$ws.Range('E16').Value = '  -3.67%  '
$ws.Range('D17').Value = '67.680.02'
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.116'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  -2.21%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '454.44'
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.46'
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.640'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.69'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').Value = '3.695.71'
$ws.Range('E25').Value = '  -3.36%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  -5.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.57'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.35'
$ws.Range('E29').Value = '  -6.46%  '
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.96'
$ws.Range('E33').Value = '  -2.52%  '
$ws.Range('E34').Value = '  -4.07%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.24'
$ws.Range('E35').Value = '  -3.14%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.158'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('D37').Value = '3.551.45'
$ws.Range('E37').Value = '  -3.16%  '
$ws.Range('E38').Value = '  -3.80%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '176.83'
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('E43').Value = '  -6.26%  '
$ws.Range('E44').Value = '  -5.50%  '
$ws.Range('E45').Value = '  -4.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.03'
$ws.Range('E46').Value = '  +6.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.99'
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('E48').Value = '  -4.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.69'
$ws.Range('E49').Value = '  -1.31%  '
$ws.Range('E50').Value = '  -5.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  -4.03%  '
